$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.209.18"
$ws.Range("E2").Value = "  +4.63%  "
$ws.Range("D3").Value = "1.905.62"
$ws.Range("E3").Value = "  +1.23%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "326.12"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.69%  "
$ws.Range("E6").Value = "  -0.30%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5160"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.36%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4017"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.47%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08460"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.93%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.77"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.13%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.120"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.81%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "23.35"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +14.27%  "
$ws.Range("E13").Value = "  +3.49%  "
$ws.Range("D14").Value = "1.917.15"
$ws.Range("E14").Value = "  +2.49%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.359"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.11%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.002"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.26%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "95.03"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.98%  "
$ws.Range("E18").Value = "  +1.12%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06712"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.00%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.36"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.82%  "
$ws.Range("E21").Value = "  -0.34%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.997"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.81%  "
$ws.Range("D23").Value = "30.220.85"
$ws.Range("E23").Value = "  +4.74%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.31"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.76%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.214"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.45%  "
$ws.Range("D26").Value = "2.134.58"
$ws.Range("E26").Value = "  +2.40%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.79"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.03%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "162.44"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.40%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.392"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.08%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "129.73"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.40%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.105"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.05%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1062"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.22%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.992"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.31%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.636"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.87%  "
$ws.Range("E35").Value = "  +1.34%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06580"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.32%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2204"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.80%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.228"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.36%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.172"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.88%  "
$ws.Range("E40").Value = "  +6.63%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.768"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.80%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6509"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.79%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.233"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.51%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6127"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.16%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.23"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.61%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.727"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.43%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.063"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.67%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.242"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.26%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "125.01"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.21%  "
$ws.Range("E50").Value = "  -0.53%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "79.31"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.67%  "
